$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1) - style matches the other header cells (bold/border/centered)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-29
$data = @(
    @(2, 7, 7),
    @(3, 7, 8),
    @(4, 8, 8),
    @(5, 10, 10),
    @(6, 7, 8),
    @(7, 7, 7),
    @(8, 7, 8),
    @(9, 8, 8),
    @(10, 7, 7),
    @(11, 7, 7),
    @(12, 9, 9),
    @(13, 9, 9),
    @(14, 8, 8),
    @(15, 6, 7),
    @(16, 8, 8),
    @(17, 7, 7),
    @(18, 6, 6),
    @(19, 5, 5),
    @(20, 9, 9),
    @(21, 8, 8),
    @(22, 7, 7),
    @(23, 7, 7),
    @(24, 8, 8),
    @(25, 6, 6),
    @(26, 7, 7),
    @(27, 8, 8),
    @(28, 4, 4),
    @(29, 3, 3)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 9).Value = $entry[1]
    $ws.Cells.Item($row, 10).Value = $entry[2]
}
